$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Row 22: H2O/CO2 column S22 cleared (no longer has a value)
$ws.Range("S22").ClearContents()

# Rows 28-31: fix P (H2O) / Q (CO2) columns for the final four analyses
$ws.Range("P28").Value = 0.5
$ws.Range("Q28").Value = 0

$ws.Range("P29").Value = 1
$ws.Range("Q29").Value = 0

$ws.Range("P30").Value = 1.5
$ws.Range("Q30").Value = 0

$ws.Range("P31").Value = 2
$ws.Range("Q31").Value = 0

# Update the saved selection to match final editing position
$ws.Range("A28").Select()
